$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 9
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 4
$ws.Range("C7").Value = 6
$ws.Range("C8").Value = 3
$ws.Range("C9").Value = 5

# Row 10: B10 text change and C10 value change
$ws.Range("B10").Value = "<senten>"
$ws.Range("C10").Value = 8

$ws.Range("C11").Value = 4
$ws.Range("C13").Value = 7
$ws.Range("C14").Value = 5
$ws.Range("C15").Value = 3
$ws.Range("C17").Value = 5
$ws.Range("C18").Value = 7
